$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: new row r gets the old D/J/K/L/M/P values of row $map[r].
$map = @{2=12; 3=41; 4=16; 5=26; 6=21; 7=39; 8=6; 9=15; 10=42; 11=23; 12=8; 13=44; 14=25; 15=2; 16=22; 17=36; 18=27; 19=33; 20=20; 21=5; 22=17; 23=14; 24=29; 25=28; 26=38; 27=30; 28=31; 29=7; 30=9; 31=34; 32=32; 33=19; 34=18; 35=3; 36=40; 37=10; 38=35; 39=4; 40=37; 41=11; 42=24; 43=43; 44=13}

# Columns being permuted: D(4) J(10) K(11) L(12) M(13) P(16)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot all the "before" values first, since the permutation mixes rows
# and we must not overwrite a source row before it has been read.
$snapshot = @{}
foreach ($r in 2..44) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in 2..44) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
